$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.024.08"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.132.99"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.128.39"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("E13").Value = "  -3.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "3.657.51"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "63.812.93"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "3.135.79"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.729"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.03%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "80.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.11%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("E32").Value = "  +3.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").Value = "0.0₃0852"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("E37").Value = "  -4.17%  "
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "439.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.284"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.30%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.901.45"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +15.26%  "
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.29%  "
